$wb = $excel.ActiveWorkbook

# UsuariosRegistro sheet - E-Mail column (C2:C6) holds the canonical
# shared-string entries that also get reused (shared) by LoginData.
$wsUsuarios = $wb.Worksheets.Item("UsuariosRegistro")
$wsUsuarios.Range("C2").Value = "juan.perez+20251111_202811@test.com"
$wsUsuarios.Range("C3").Value = "maria.gonzalez+20251111_202811@test.com"
$wsUsuarios.Range("C4").Value = "carlos.rodriguez+20251111_202811@test.com"
$wsUsuarios.Range("C5").Value = "ana.martinez+20251111_202811@test.com"
$wsUsuarios.Range("C6").Value = "luis.garcia+20251111_202811@test.com"

# LoginData sheet reuses the same two e-mail strings (juan.perez and
# maria.gonzalez) in column A, rows 2 and 3 - update them to match.
$wsLogin = $wb.Worksheets.Item("LoginData")
$wsLogin.Range("A2").Value = "juan.perez+20251111_202811@test.com"
$wsLogin.Range("A3").Value = "maria.gonzalez+20251111_202811@test.com"
